$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "edu-launch-esp32s3"

$ws.Range("A6").Value = 45042
$ws.Range("A7").Value = 45042
$ws.Range("A8").Value = 45042

$ws.Range("C10").Value = "Issues test needs HAS PSAM"

$ws.Columns("A").ColumnWidth = 30.6328125

$ws.Range("C10").Select()

$wb.Windows.Item(1).DisplayZoom = 96
